$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.20211632131014
$ws.Range("C2").Value = 11.69932770511651
$ws.Range("D2").Value = 4.581185571546554
$ws.Range("F2").Value = 49.95198945433712
$ws.Range("G2").Value = 3.75782654549194
$ws.Range("I2").Value = 35.01219703964031
$ws.Range("J2").Value = 9.498795707720062
$ws.Range("L2").Value = 12.23761751769462
$ws.Range("M2").Value = 20.45175224154402
$ws.Range("N2").Value = 21.84510553373602
$ws.Range("B3").Value = 23.92119859376479
$ws.Range("C3").Value = 11.30429259650018
$ws.Range("D3").Value = 4.5498350858387
$ws.Range("F3").Value = 49.92958505205166
$ws.Range("G3").Value = 3.761661380767777
$ws.Range("I3").Value = 35.06539113271728
$ws.Range("J3").Value = 9.503238491592846
$ws.Range("L3").Value = 12.25298376601398
$ws.Range("M3").Value = 20.41213776007392
$ws.Range("N3").Value = 21.91592544245623
$ws.Range("B4").Value = 23.75341820086984
$ws.Range("C4").Value = 11.05837847861292
$ws.Range("D4").Value = 4.530011490559326
$ws.Range("F4").Value = 49.9276744132593
$ws.Range("G4").Value = 3.76413895203808
$ws.Range("I4").Value = 35.10523342111419
$ws.Range("J4").Value = 9.5061454540604
$ws.Range("L4").Value = 12.26399994906241
$ws.Range("M4").Value = 20.39171788041162
$ws.Range("N4").Value = 21.96140761406224
$ws.Range("B5").Value = 23.68630691916077
$ws.Range("C5").Value = 10.95749872650325
$ws.Range("D5").Value = 4.521789013261047
$ws.Range("F5").Value = 49.92987167662958
$ws.Range("G5").Value = 3.765179619031572
$ws.Range("I5").Value = 35.12326980830006
$ws.Range("J5").Value = 9.507375189532384
$ws.Range("L5").Value = 12.26888703571757
$ws.Range("M5").Value = 20.38438373989148
$ws.Range("N5").Value = 21.98044581986823
$ws.Range("B6").Value = 23.67524149625504
$ws.Range("C6").Value = 10.94071283067325
$ws.Range("D6").Value = 4.520414973656799
$ws.Range("F6").Value = 49.93041614456312
$ws.Range("G6").Value = 3.765354298920001
$ws.Range("I6").Value = 35.126373306933
$ws.Range("J6").Value = 9.507582114389095
$ws.Range("L6").Value = 12.26972257261641
$ws.Range("M6").Value = 20.3832256781451
$ws.Range("N6").Value = 21.98363756978153
$ws.Range("B7").Value = 23.75250790999612
$ws.Range("C7").Value = 11.05702043859238
$ws.Range("D7").Value = 4.529901182922416
$ws.Range("F7").Value = 49.92769200255678
$ws.Range("G7").Value = 3.76415286101924
$ws.Range("I7").Value = 35.10546938422574
$ws.Range("J7").Value = 9.506161855865997
$ws.Range("L7").Value = 12.2640642465888
$ws.Range("M7").Value = 20.39161496609811
$ws.Range("N7").Value = 21.96166232819175
$ws.Range("B8").Value = 24.10432689453631
$ws.Range("C8").Value = 11.56391946679366
$ws.Range("D8").Value = 4.570495231532203
$ws.Range("F8").Value = 49.9418046715688
$ws.Range("G8").Value = 3.75912334460274
$ws.Range("I8").Value = 35.02904486202659
$ws.Range("J8").Value = 9.500290469442255
$ws.Range("L8").Value = 12.24258779495523
$ws.Range("M8").Value = 20.43728653114304
$ws.Range("N8").Value = 21.86911025729421
$ws.Range("B9").Value = 24.82804413884364
$ws.Range("C9").Value = 12.52417524770628
$ws.Range("D9").Value = 4.645536664255124
$ws.Range("F9").Value = 50.06352048880671
$ws.Range("G9").Value = 3.750230897533382
$ws.Range("I9").Value = 34.93637864347546
$ws.Range("J9").Value = 9.4901933313727
$ws.Range("L9").Value = 12.21300583158247
$ws.Range("M9").Value = 20.55754644148867
$ws.Range("N9").Value = 21.7034111348968
$ws.Range("B10").Value = 25.37553470673668
$ws.Range("C10").Value = 13.20056392935909
$ws.Range("D10").Value = 4.697881901930423
$ws.Range("F10").Value = 50.21023154226495
$ws.Range("G10").Value = 3.744281842343975
$ws.Range("I10").Value = 34.90345356489681
$ws.Range("J10").Value = 9.483632743780651
$ws.Range("L10").Value = 12.19889393437831
$ws.Range("M10").Value = 20.66420207963118
$ws.Range("N10").Value = 21.59121293820062
$ws.Range("B11").Value = 25.62698709318564
$ws.Range("C11").Value = 13.50042756915322
$ws.Range("D11").Value = 4.721088754711127
$ws.Range("F11").Value = 50.28936256452801
$ws.Range("G11").Value = 3.741700749226765
$ws.Range("I11").Value = 34.89616034374979
$ws.Range("J11").Value = 9.48083319880392
$ws.Range("L11").Value = 12.19412428747583
$ws.Range("M11").Value = 20.71659268266569
$ws.Range("N11").Value = 21.54222405385378
$ws.Range("B12").Value = 25.72246334485488
$ws.Range("C12").Value = 13.61274359957669
$ws.Range("D12").Value = 4.729789374201364
$ws.Range("F12").Value = 50.32110132734814
$ws.Range("G12").Value = 3.740741233780136
$ws.Range("I12").Value = 34.89450686266701
$ws.Range("J12").Value = 9.479799580299773
$ws.Range("L12").Value = 12.19255489540827
$ws.Range("M12").Value = 20.73697849924459
$ws.Range("N12").Value = 21.52396676620429
$ws.Range("B13").Value = 25.70189079680702
$ws.Range("C13").Value = 13.58861091512919
$ws.Range("D13").Value = 4.727919427698444
$ws.Range("F13").Value = 50.31418709160305
$ws.Range("G13").Value = 3.740947088743755
$ws.Range("I13").Value = 34.89481363102582
$ws.Range("J13").Value = 9.480021010827571
$ws.Range("L13").Value = 12.19288236999495
$ws.Range("M13").Value = 20.73256390734848
$ws.Range("N13").Value = 21.52788575354831
$ws.Range("B14").Value = 25.63483733689341
$ws.Range("C14").Value = 13.50969322989266
$ws.Range("D14").Value = 4.721806311337924
$ws.Range("F14").Value = 50.29193824181433
$ws.Range("G14").Value = 3.741621451427265
$ws.Range("I14").Value = 34.89600208203484
$ws.Range("J14").Value = 9.480747631571818
$ws.Range("L14").Value = 12.19399043040788
$ws.Range("M14").Value = 20.7182589417405
$ws.Range("N14").Value = 21.5407161352971
$ws.Range("B15").Value = 25.59379599535097
$ws.Range("C15").Value = 13.46118987780071
$ws.Range("D15").Value = 4.718050472344192
$ws.Range("F15").Value = 50.27854086690353
$ws.Range("G15").Value = 3.742036844408357
$ws.Range("I15").Value = 34.89687446413798
$ws.Range("J15").Value = 9.481196157590025
$ws.Range("L15").Value = 12.1946999684438
$ws.Range("M15").Value = 20.70956760976465
$ws.Range("N15").Value = 21.54861333464309
$ws.Range("B16").Value = 25.35914238449859
$ws.Range("C16").Value = 13.18079990164013
$ws.Range("D16").Value = 4.696353076419058
$ws.Range("F16").Value = 50.20530882821165
$ws.Range("G16").Value = 3.744453032044281
$ws.Range("I16").Value = 34.90408511937592
$ws.Range("J16").Value = 9.48381941441696
$ws.Range("L16").Value = 12.19923880225564
$ws.Range("M16").Value = 20.66085529616837
$ws.Range("N16").Value = 21.59445564324528
$ws.Range("B17").Value = 25.21574089409628
$ws.Range("C17").Value = 13.00670221238516
$ws.Range("D17").Value = 4.682887207284249
$ws.Range("F17").Value = 50.16355183509484
$ws.Range("G17").Value = 3.745967263702904
$ws.Range("I17").Value = 34.91047932564496
$ws.Range("J17").Value = 9.485475997393015
$ws.Range("L17").Value = 12.20244548496057
$ws.Range("M17").Value = 20.63195664405392
$ws.Range("N17").Value = 21.62310281560999
$ws.Range("B18").Value = 25.13349090403078
$ws.Range("C18").Value = 12.90583467168626
$ws.Range("D18").Value = 4.675085041995989
$ws.Range("F18").Value = 50.14070128753647
$ws.Range("G18").Value = 3.746849997061704
$ws.Range("I18").Value = 34.91488025024069
$ws.Range("J18").Value = 9.486446228391779
$ws.Range("L18").Value = 12.20444521537223
$ws.Range("M18").Value = 20.61570005122595
$ws.Range("N18").Value = 21.63977299428341
$ws.Range("B19").Value = 25.1056847962063
$ws.Range("C19").Value = 12.87156064002465
$ws.Range("D19").Value = 4.672433589595902
$ws.Range("F19").Value = 50.13316513198387
$ws.Range("G19").Value = 3.747150903247058
$ws.Range("I19").Value = 34.91649441995175
$ws.Range("J19").Value = 9.486777724281149
$ws.Range("L19").Value = 12.20514898202352
$ws.Range("M19").Value = 20.61025886235653
$ws.Range("N19").Value = 21.64545042748885
$ws.Range("B20").Value = 25.23098298530428
$ws.Range("C20").Value = 13.02531171740883
$ws.Range("D20").Value = 4.68432656847471
$ws.Range("F20").Value = 50.16787621529415
$ws.Range("G20").Value = 3.745804852010914
$ws.Range("I20").Value = 34.9097237834316
$ws.Range("J20").Value = 9.485297850197096
$ws.Range("L20").Value = 12.20208805530444
$ws.Range("M20").Value = 20.63499523540956
$ws.Range("N20").Value = 21.62003330030964
$ws.Range("B21").Value = 25.65452629189362
$ws.Range("C21").Value = 13.53290762153438
$ws.Range("D21").Value = 4.723604253598501
$ws.Range("F21").Value = 50.29842521086545
$ws.Range("G21").Value = 3.741422890142504
$ws.Range("I21").Value = 34.89562290276817
$ws.Range("J21").Value = 9.480533486644859
$ws.Range("L21").Value = 12.19365854448897
$ws.Range("M21").Value = 20.72244590633008
$ws.Range("N21").Value = 21.53693957814199
$ws.Range("B22").Value = 25.93279263236528
$ws.Range("C22").Value = 13.85740522123229
$ws.Range("D22").Value = 4.748765446224442
$ws.Range("F22").Value = 50.39407888818474
$ws.Range("G22").Value = 3.738663242091002
$ws.Range("I22").Value = 34.89286832979144
$ws.Range("J22").Value = 9.477574174895906
$ws.Range("L22").Value = 12.1895291835736
$ws.Range("M22").Value = 20.7827809423243
$ws.Range("N22").Value = 21.48434465100802
$ws.Range("B23").Value = 25.78417212251613
$ws.Range("C23").Value = 13.68491086774708
$ws.Range("D23").Value = 4.735383111730055
$ws.Range("F23").Value = 50.34208458907192
$ws.Range("G23").Value = 3.740126618392555
$ws.Range("I23").Value = 34.89374636150811
$ws.Range("J23").Value = 9.479139506813009
$ws.Range("L23").Value = 12.19160702297849
$ws.Range("M23").Value = 20.75029149301811
$ws.Range("N23").Value = 21.51225931146248
$ws.Range("B24").Value = 25.2240914266027
$ws.Range("C24").Value = 13.01690077025274
$ws.Range("D24").Value = 4.683676021602895
$ws.Range("F24").Value = 50.16591756118311
$ws.Range("G24").Value = 3.745878240384209
$ws.Range("I24").Value = 34.91006310654465
$ws.Range("J24").Value = 9.485378334973653
$ws.Range("L24").Value = 12.20224916278053
$ws.Range("M24").Value = 20.63362037350799
$ws.Range("N24").Value = 21.62142040339494
$ws.Range("B25").Value = 24.6291357508983
$ws.Range("C25").Value = 12.26895078321884
$ws.Range("D25").Value = 4.625725097103915
$ws.Range("F25").Value = 50.02052180507101
$ws.Range("G25").Value = 3.752533413565568
$ws.Range("I25").Value = 34.95529209198155
$ws.Range("J25").Value = 9.492773822634986
$ws.Range("L25").Value = 12.219668498878
$ws.Range("M25").Value = 20.55754644148867
$ws.Range("N25").Value = 21.74655533630508
